$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 104, shifting existing rows 104:137 down to 107:140
$ws.Range("A104:T106").EntireRow.Insert()

# Populate the new rows 104:106 with this week's new data (date 2021-11-23 / serial 44523)
$ws.Range("A104").Value = 3
$ws.Range("B104").Value = "Femacal de La Calera"
$ws.Range("C104").Value = "Coquimbo"
$ws.Range("D104").Value = 44523
$ws.Range("E104").Value = 5
$ws.Range("F104").Value = "Fruta"
$ws.Range("G104").Value = 100107
$ws.Range("H104").Value = "Otros"
$ws.Range("I104").Value = 100107002
$ws.Range("J104").Value = "Chirimoya"
$ws.Range("K104").Value = "Cultivar IV Región"
$ws.Range("L104").Value = "Especial"
$ws.Range("M104").Value = 80
$ws.Range("N104").Value = 26000
$ws.Range("O104").Value = 26000
$ws.Range("P104").Value = 26000
$ws.Range("Q104").Value = "$/bandeja 10 kilos"
$ws.Range("R104").Value = "Provincia de Limarí"
$ws.Range("S104").Value = 2600
$ws.Range("T104").Value = 10

$ws.Range("A105").Value = 3
$ws.Range("B105").Value = "Femacal de La Calera"
$ws.Range("C105").Value = "Coquimbo"
$ws.Range("D105").Value = 44523
$ws.Range("E105").Value = 5
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100107
$ws.Range("H105").Value = "Otros"
$ws.Range("I105").Value = 100107002
$ws.Range("J105").Value = "Chirimoya"
$ws.Range("K105").Value = "Cultivar IV Región"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 95
$ws.Range("N105").Value = 24000
$ws.Range("O105").Value = 24000
$ws.Range("P105").Value = 24000
$ws.Range("Q105").Value = "$/bandeja 10 kilos"
$ws.Range("R105").Value = "Provincia de Limarí"
$ws.Range("S105").Value = 2400
$ws.Range("T105").Value = 10

$ws.Range("A106").Value = 3
$ws.Range("B106").Value = "Femacal de La Calera"
$ws.Range("C106").Value = "Coquimbo"
$ws.Range("D106").Value = 44523
$ws.Range("E106").Value = 5
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100107
$ws.Range("H106").Value = "Otros"
$ws.Range("I106").Value = 100107002
$ws.Range("J106").Value = "Chirimoya"
$ws.Range("K106").Value = "Cultivar IV Región"
$ws.Range("L106").Value = "Segunda"
$ws.Range("M106").Value = 95
$ws.Range("N106").Value = 21000
$ws.Range("O106").Value = 21000
$ws.Range("P106").Value = 21000
$ws.Range("Q106").Value = "$/bandeja 10 kilos"
$ws.Range("R106").Value = "Provincia de Limarí"
$ws.Range("S106").Value = 2100
$ws.Range("T106").Value = 10
